$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text (matching the source data's
# inline-string representation) instead of Excel auto-coercing numeric-looking
# strings ("299.17") or percentages ("-1.14%") into numbers.
$dataRange = $ws.Range("D2:E50")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "299.17"
$ws.Range("E2").Value = "-1.14%"
$ws.Range("D3").Value = "31.47"
$ws.Range("E3").Value = "-1.10%"
$ws.Range("D4").Value = "5.091"
$ws.Range("E4").Value = "-1.42%"
$ws.Range("D5").Value = "0.07856"
$ws.Range("E5").Value = "0.13%"
$ws.Range("D6").Value = "2.273"
$ws.Range("E6").Value = "-1.69%"
$ws.Range("D7").Value = "7.809"
$ws.Range("E7").Value = "-1.71%"
$ws.Range("D8").Value = "3.853"
$ws.Range("E8").Value = "-0.41%"
$ws.Range("D9").Value = "0.9234"
$ws.Range("E9").Value = "1.73%"
$ws.Range("D10").Value = "0.1747"
$ws.Range("E10").Value = "0.30%"
$ws.Range("D11").Value = "0.07603"
$ws.Range("E11").Value = "3.50%"
$ws.Range("D12").Value = "0.09236"
$ws.Range("E12").Value = "13.22%"
$ws.Range("D13").Value = "0.03003"
$ws.Range("E13").Value = "-3.80%"
$ws.Range("E14").Value = "0.68%"
$ws.Range("D15").Value = "0.001508"
$ws.Range("E15").Value = "-0.14%"
$ws.Range("D16").Value = "0.006011"
$ws.Range("E16").Value = "0.25%"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").Value = "-0.72%"
$ws.Range("D18").Value = "2.244"
$ws.Range("E18").Value = "0.08%"
$ws.Range("E20").Value = "-2.13%"
$ws.Range("D21").Value = "3.912"
$ws.Range("E21").Value = "-16.30%"
$ws.Range("D23").Value = "0.04620"
$ws.Range("E23").Value = "-0.80%"
$ws.Range("E24").Value = "-0.91%"
$ws.Range("D25").Value = "0.004466"
$ws.Range("E25").Value = "-1.18%"
$ws.Range("E26").Value = "-7.30%"
$ws.Range("D27").Value = "0.0003399"
$ws.Range("E27").Value = "23.95%"
$ws.Range("D39").Value = "0.01737"
$ws.Range("E39").Value = "-5.85%"
$ws.Range("D40").Value = "0.04609"
$ws.Range("E40").Value = "0.83%"
$ws.Range("D41").Value = "0.007041"
$ws.Range("E41").Value = "-3.87%"
$ws.Range("E42").Value = "-0.27%"
$ws.Range("D43").Value = "0.002191"
$ws.Range("E43").Value = "-2.60%"
$ws.Range("D44").Value = "0.009743"
$ws.Range("E44").Value = "-9.28%"
$ws.Range("D45").Value = "0.00006291"
$ws.Range("E45").Value = "-2.56%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "0.05%"
$ws.Range("D47").Value = "0.007987"
$ws.Range("E47").Value = "-19.30%"
$ws.Range("D48").Value = "1.154"
$ws.Range("E48").Value = "40.64%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.05%"

# Reset the style back to the default (Normal) so only the cell values changed -
# matches the original workbook, which had no explicit number format on this range.
$dataRange.Style = "Normal"
